# add save column in s_vals sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new
# "Save" header (H1) so it keeps the same bold/centered/bordered style,
# then overwrite the copied text with the new header label.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for the Save column.
$ws.Range("H2").Value = 1
